# Scheduled data refresh: update market/profit figures (columns H-N) on the
# Leve profitability sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 193
$ws.Range("I2").Value = 86
$ws.Range("K2").Value = 86
$ws.Range("M2").Value = 27
$ws.Range("H18").Value = 592.7143
$ws.Range("I18").Value = 592.7143
$ws.Range("K18").Value = 592.7143
$ws.Range("M18").Value = -308.7143
$ws.Range("H32").Value = 1221.0555
$ws.Range("I32").Value = 497
$ws.Range("J32").Value = 1681.8182
$ws.Range("K32").Value = 497
$ws.Range("L32").Value = 1681.8182
$ws.Range("M32").Value = -171
$ws.Range("N32").Value = -2333.8182
$ws.Range("H116").Value = 14350.625
$ws.Range("I116").Value = 26501.25
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 26501.25
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = -23059.25
$ws.Range("N116").Value = -9084
$ws.Range("H137").Value = 1223.0555
$ws.Range("I137").Value = 1028.4138
$ws.Range("J137").Value = 2029.4286
$ws.Range("K137").Value = 3085.2414
$ws.Range("L137").Value = 6088.2858
$ws.Range("M137").Value = -535.2413999999999
$ws.Range("N137").Value = -11188.2858
$ws.Range("H138").Value = 4510.82
$ws.Range("I138").Value = 973.8214
$ws.Range("J138").Value = 9012.454
$ws.Range("K138").Value = 2921.4642
$ws.Range("L138").Value = 27037.362
$ws.Range("M138").Value = 2218.5358
$ws.Range("N138").Value = -37317.362
$ws.Range("H141").Value = 1840.5588
$ws.Range("I141").Value = 1563.1786
$ws.Range("J141").Value = 3135
$ws.Range("K141").Value = 4689.5358
$ws.Range("L141").Value = 9405
$ws.Range("M141").Value = 490.4642000000003
$ws.Range("N141").Value = -19765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4775.4062
$ws.Range("I32").Value = 3580.5818
$ws.Range("J32").Value = 12077.111
$ws.Range("K32").Value = 3580.5818
$ws.Range("L32").Value = 12077.111
$ws.Range("M32").Value = -3293.5818
$ws.Range("N32").Value = -12651.111
$ws.Range("H61").Value = 7982.1177
$ws.Range("I61").Value = 9275.286
$ws.Range("K61").Value = 9275.286
$ws.Range("M61").Value = -9063.286
$ws.Range("H74").Value = 1178
$ws.Range("I74").Value = 1151.5
$ws.Range("K74").Value = 1151.5
$ws.Range("M74").Value = -277.5
$ws.Range("H77").Value = 1178
$ws.Range("I77").Value = 1151.5
$ws.Range("K77").Value = 5757.5
$ws.Range("M77").Value = -1389.5
$ws.Range("H132").Value = 6672.125
$ws.Range("I132").Value = 5462.4
$ws.Range("K132").Value = 16387.2
$ws.Range("M132").Value = -13857.2
$ws.Range("H136").Value = 7982.1177
$ws.Range("I136").Value = 9275.286
$ws.Range("K136").Value = 27825.858
$ws.Range("M136").Value = -25275.858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5476.0645
$ws.Range("I134").Value = 7686.6113
$ws.Range("K134").Value = 23059.8339
$ws.Range("M134").Value = -20524.8339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 731.25
$ws.Range("I22").Value = 150
$ws.Range("J22").Value = 814.2857
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 814.2857
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = -1514.2857
$ws.Range("H31").Value = 5165.684
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5165.684
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5165.684
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -5755.684
$ws.Range("H34").Value = 5165.684
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5165.684
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5165.684
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -5569.684
$ws.Range("H58").Value = 1209.875
$ws.Range("I58").Value = 987.56525
$ws.Range("K58").Value = 987.56525
$ws.Range("M58").Value = -784.56525
$ws.Range("H132").Value = 2770.4814
$ws.Range("I132").Value = 2633.6316
$ws.Range("J132").Value = 3095.5
$ws.Range("K132").Value = 7900.8948
$ws.Range("L132").Value = 9286.5
$ws.Range("M132").Value = -5370.8948
$ws.Range("N132").Value = -14346.5
$ws.Range("H134").Value = 4728.357
$ws.Range("I134").Value = 5788.3
$ws.Range("J134").Value = 2078.5
$ws.Range("K134").Value = 17364.9
$ws.Range("L134").Value = 6235.5
$ws.Range("M134").Value = -14829.9
$ws.Range("N134").Value = -11305.5
$ws.Range("H136").Value = 1209.875
$ws.Range("I136").Value = 987.56525
$ws.Range("K136").Value = 2962.69575
$ws.Range("M136").Value = -412.6957499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 931.78
$ws.Range("I131").Value = 626.94116
$ws.Range("J131").Value = 994.21686
$ws.Range("K131").Value = 1880.82348
$ws.Range("L131").Value = 2982.65058
$ws.Range("M131").Value = 3159.17652
$ws.Range("N131").Value = -13062.65058
$ws.Range("H140").Value = 6645.273
$ws.Range("I140").Value = 3309.8
$ws.Range("K140").Value = 9929.400000000001
$ws.Range("M140").Value = -4749.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5806.6
$ws.Range("I70").Value = 5721.4165
$ws.Range("J70").Value = 5992.4546
$ws.Range("K70").Value = 5721.4165
$ws.Range("L70").Value = 5992.4546
$ws.Range("M70").Value = -5451.4165
$ws.Range("N70").Value = -6532.4546
$ws.Range("H73").Value = 5806.6
$ws.Range("I73").Value = 5721.4165
$ws.Range("J73").Value = 5992.4546
$ws.Range("K73").Value = 5721.4165
$ws.Range("L73").Value = 5992.4546
$ws.Range("M73").Value = -4785.4165
$ws.Range("N73").Value = -7864.4546
$ws.Range("H132").Value = 3210.2856
$ws.Range("I132").Value = 4095.111
$ws.Range("J132").Value = 2546.6667
$ws.Range("K132").Value = 12285.333
$ws.Range("L132").Value = 7640.000100000001
$ws.Range("M132").Value = -9755.332999999999
$ws.Range("N132").Value = -12700.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 524.9545
$ws.Range("I16").Value = 504.2353
$ws.Range("J16").Value = 595.4
$ws.Range("K16").Value = 504.2353
$ws.Range("L16").Value = 595.4
$ws.Range("M16").Value = -334.2353
$ws.Range("N16").Value = -935.4
$ws.Range("H132").Value = 9552772
$ws.Range("I132").Value = 15810048
$ws.Range("K132").Value = 47430144
$ws.Range("M132").Value = -47427614
$ws.Range("H136").Value = 6308.737
$ws.Range("I136").Value = 6487.409
$ws.Range("J136").Value = 6063.0625
$ws.Range("K136").Value = 19462.227
$ws.Range("L136").Value = 18189.1875
$ws.Range("M136").Value = -16912.227
$ws.Range("N136").Value = -23289.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1725.4872
$ws.Range("I132").Value = 1049.2778
$ws.Range("K132").Value = 3147.8334
$ws.Range("M132").Value = -617.8334000000004
$ws.Range("H136").Value = 2497.725
$ws.Range("I136").Value = 3251.5264
$ws.Range("J136").Value = 1815.7142
$ws.Range("K136").Value = 9754.5792
$ws.Range("L136").Value = 5447.142599999999
$ws.Range("M136").Value = -7204.5792
$ws.Range("N136").Value = -10547.1426
